$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8: "Good Morning" -> "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Update the selection to E8 (matches the sheetView <selection activeCell="E8" sqref="E8"/>)
$ws.Range("E8").Select()
